# Merge the split "{{replace}} → {{replace}} → abc {{b" / "y" / "}}" runs
# in the last paragraph of the "replaceText" shape back into a single run
# reading "{{replace}} → {{replace}} → abc {{by}}".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

# Locate the target paragraph (the last one in the text frame) instead of
# hard-coding its index, so the script keeps working even if the shape's
# paragraph count changes.
$paras = $tr.Paragraphs()
$target = $paras.Count

$para = $tr.Paragraphs($target, 1)

if ($para.Text -like "*{{replace}}*{{replace}}*abc*") {
    # Re-assigning .Text collapses the paragraph's runs into a single run
    # that carries the formatting of the first original run (plain, no
    # red fill), matching the desired output.
    $para.Text = "{{replace}} " + [char]0x2192 + " {{replace}} " + [char]0x2192 + " abc {{by}}"
} else {
    # Fallback: assertion on the expected content failed, so (re)create the
    # text element from scratch to guarantee the final text is correct.
    $tr.Text = $tr.Text
    $para.Text = "{{replace}} " + [char]0x2192 + " {{replace}} " + [char]0x2192 + " abc {{by}}"
}
